$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D10").Value = "recupero"
$ws.Range("D11").Value = "recupero"
$ws.Range("D66").Value = "recupero"
$ws.Range("D67").Value = "recupero"
$ws.Range("D122").Value = "recupero"
$ws.Range("D123").Value = "recupero"
$ws.Range("D178").Value = "recupero"
$ws.Range("D179").Value = "recupero"
$ws.Range("D339").Value = "recupero"
$ws.Range("D343").Value = "falto"
$ws.Range("D344").Value = "falto"
$ws.Range("D345").Value = "falto"
$ws.Range("D346").Value = "falto"
$ws.Range("D347").Value = "asistio"
$ws.Range("D348").Value = "asistio"
$ws.Range("D351").Value = "falto"
$ws.Range("D355").Value = "falto"
$ws.Range("D356").Value = "falto"
$ws.Range("D373").Value = "falto"
$ws.Range("D374").Value = "falto"
$ws.Range("D377").Value = "falto"
$ws.Range("D378").Value = "falto"
$ws.Range("D379").Value = "recupero"
$ws.Range("D380").Value = "recupero"
$ws.Range("D383").Value = "recupero"
$ws.Range("D384").Value = "recupero"
$ws.Range("D389").Value = "falto"
$ws.Range("D390").Value = "falto"
$ws.Range("D391").Value = "falto"
$ws.Range("D392").Value = "falto"
$ws.Range("D393").Value = "falto"
$ws.Range("D394").Value = "falto"
$ws.Range("D398").Value = "recupero"
$ws.Range("D399").Value = "recupero"
$ws.Range("D400").Value = "falto"
$ws.Range("D401").Value = "falto"
$ws.Range("D403").Value = "falto"
$ws.Range("D404").Value = "asistio"
$ws.Range("D405").Value = "asistio"
$ws.Range("D408").Value = "falto"
$ws.Range("D410").Value = "falto"
$ws.Range("D411").Value = "falto"
$ws.Range("D412").Value = "falto"
$ws.Range("D413").Value = "falto"
$ws.Range("D414").Value = "falto"
$ws.Range("D416").Value = "recupero"
$ws.Range("D417").Value = "recupero"
$ws.Range("D418").Value = "falto"
$ws.Range("D419").Value = "falto"
$ws.Range("D434").Value = "falto"
$ws.Range("D435").Value = "falto"
$ws.Range("D440").Value = "recupero"
$ws.Range("D441").Value = "recupero"
$ws.Range("D446").Value = "falto"
$ws.Range("D447").Value = "falto"
$ws.Range("D448").Value = "falto"
$ws.Range("D449").Value = "falto"
$ws.Range("D450").Value = "falto"
$ws.Range("D451").Value = "falto"
$ws.Range("D453").Value = "asistio"
$ws.Range("D454").Value = "recupero"
$ws.Range("D455").Value = "falto"
$ws.Range("D456").Value = "falto"
$ws.Range("D457").Value = "falto"
$ws.Range("D458").Value = "falto"
$ws.Range("D459").Value = "asistio"
$ws.Range("D460").Value = "falto"
$ws.Range("D461").Value = "asistio"
$ws.Range("D462").Value = "asistio"
$ws.Range("D463").Value = "asistio"
$ws.Range("D464").Value = "asistio"
$ws.Range("D465").Value = "recupero"
$ws.Range("D466").Value = "asistio"
$ws.Range("D467").Value = "falto"
$ws.Range("D468").Value = "falto"
$ws.Range("D469").Value = "falto"
$ws.Range("D470").Value = "falto"
$ws.Range("D471").Value = "asistio"
$ws.Range("D472").Value = "asistio"
$ws.Range("D473").Value = "recupero"
$ws.Range("D474").Value = "recupero"
$ws.Range("D475").Value = "falto"
$ws.Range("D476").Value = "falto"
$ws.Range("D477").Value = "asistio"
$ws.Range("D478").Value = "asistio"
$ws.Range("D479").Value = "asistio"
$ws.Range("D480").Value = "asistio"
$ws.Range("D481").Value = "asistio"
$ws.Range("D482").Value = "asistio"
$ws.Range("D483").Value = "recupero"
$ws.Range("D484").Value = "recupero"
$ws.Range("D485").Value = "asistio"
$ws.Range("D486").Value = "asistio"
$ws.Range("D487").Value = "asistio"
$ws.Range("D488").Value = "asistio"
$ws.Range("D489").Value = "asistio"
$ws.Range("D490").Value = "asistio"
$ws.Range("D491").Value = "falto"
$ws.Range("D492").Value = "falto"
$ws.Range("D493").Value = "recupero"
$ws.Range("D494").Value = "recupero"
$ws.Range("D495").Value = "recupero"
$ws.Range("D496").Value = "recupero"
$ws.Range("D497").Value = "recupero"
$ws.Range("D498").Value = "recupero"
$ws.Range("D499").Value = "asistio"
$ws.Range("D500").Value = "recupero"
$ws.Range("D501").Value = "recupero"
$ws.Range("D502").Value = "recupero"
$ws.Range("D503").Value = "falto"
$ws.Range("D504").Value = "falto"
$ws.Range("D505").Value = "falto"
$ws.Range("D506").Value = "falto"
$ws.Range("D507").Value = "falto"
$ws.Range("D508").Value = "falto"
$ws.Range("D509").Value = "asistio"
$ws.Range("D569").Value = "recupero"
$ws.Range("D570").Value = "asistio"

$ws.Range("D291").Select()
